# AssignmentSecondaryStructure.docx edit
#
# 1. Fix the typos "pars" -> "pairs" and "occurrence" -> "occurrences" in the
#    "...all other pars have zero occurrence. " sentence, reproducing the
#    run fragmentation that Word leaves behind when the two words were
#    corrected interactively (several small runs instead of one long run).
# 2. Move the "_GoBack" bookmark (Word's "last edit location" marker) from
#    the end of the "...by choosing GPU." paragraph to sit right after the
#    newly-corrected sentence, between "occurrences" and the final ". ".

$d = $word.ActiveDocument

# --- Step 1: correct the wording (this naturally merges into one run) ---
$rng = $d.Content
$rng.Find.Execute("ll other pars have zero occurrence. ", $true, $false, `
    $false, $false, $false, $true, 1, $false, `
    "ll other pairs have zero occurrences. ", 2)

# --- Step 2: locate the corrected sentence again to get fresh offsets ---
$rng2 = $d.Content
$rng2.Find.Execute("all other pairs have zero occurrences. ", $true, $false, `
    $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng2.Start

# --- Step 3: re-split the merged run back into the expected fragments ---
# "a" | "ll other pa" | "i" | "rs have zero occurrence" | "s" | ". "
# Temporary bookmarks force a run boundary at each offset; once all the
# boundaries we need exist we discard the temporary ones (removing a
# bookmark does not re-merge the surrounding runs).
$d.Bookmarks.Add("ZZTmpSplit0", $d.Range($start + 1, $start + 1))
$d.Bookmarks.Add("ZZTmpSplit1", $d.Range($start + 12, $start + 12))
$d.Bookmarks.Add("ZZTmpSplit2", $d.Range($start + 13, $start + 13))
$d.Bookmarks.Add("ZZTmpSplit3", $d.Range($start + 36, $start + 36))
$d.Bookmarks.Add("ZZTmpSplit4", $d.Range($start + 37, $start + 37))

# Remember where the last split point (after "s", before ". ") is -- that
# is where "_GoBack" belongs.
$goBackRange = $d.Bookmarks("ZZTmpSplit4").Range

$d.Bookmarks("ZZTmpSplit0").Delete()
$d.Bookmarks("ZZTmpSplit1").Delete()
$d.Bookmarks("ZZTmpSplit2").Delete()
$d.Bookmarks("ZZTmpSplit3").Delete()
$d.Bookmarks("ZZTmpSplit4").Delete()

# --- Step 4: (re)create "_GoBack" at the new location ---
# A document only ever has a single "_GoBack" bookmark, so adding it here
# also removes the old one that used to sit after "...by choosing GPU."
$d.Bookmarks.Add("_GoBack", $goBackRange)
